$wb = $excel.ActiveWorkbook

# This script applies cached-value updates to the profit-tracking sheets,
# mirroring a scheduled data refresh (no formulas are present in this
# workbook; every cell below is a literal number produced by an external
# pricing job).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1723.5714
$ws.Range("J19").Value = 1713.5
$ws.Range("L19").Value = 1713.5
$ws.Range("N19").Value = -2063.5
$ws.Range("H40").Value = 3747.4443
$ws.Range("J40").Value = 4023.9565
$ws.Range("L40").Value = 4023.9565
$ws.Range("N40").Value = -4373.9565
$ws.Range("H51").Value = 4714.4443
$ws.Range("J51").Value = 4395
$ws.Range("L51").Value = 4395
$ws.Range("N51").Value = -5363
$ws.Range("H64").Value = 9925.385
$ws.Range("I64").Value = 5005
$ws.Range("K64").Value = 5005
$ws.Range("M64").Value = -4757
$ws.Range("H67").Value = 9925.385
$ws.Range("I67").Value = 5005
$ws.Range("K67").Value = 5005
$ws.Range("M67").Value = -4147
$ws.Range("H87").Value = 60294.5
$ws.Range("J87").Value = 60294.5
$ws.Range("L87").Value = 60294.5
$ws.Range("N87").Value = -62790.5
$ws.Range("H90").Value = 60294.5
$ws.Range("J90").Value = 60294.5
$ws.Range("L90").Value = 180883.5
$ws.Range("N90").Value = -193363.5
$ws.Range("H132").Value = 2216.9546
$ws.Range("I132").Value = 2083.842
$ws.Range("J132").Value = 3060
$ws.Range("K132").Value = 6251.526
$ws.Range("L132").Value = 9180
$ws.Range("M132").Value = -3721.526
$ws.Range("N132").Value = -14240
$ws.Range("H137").Value = 1307.2307
$ws.Range("I137").Value = 1370
$ws.Range("J137").Value = 962
$ws.Range("K137").Value = 4110
$ws.Range("L137").Value = 2886
$ws.Range("M137").Value = -1560
$ws.Range("N137").Value = -7986
$ws.Range("H138").Value = 2726.8572
$ws.Range("I138").Value = 1835.8667
$ws.Range("J138").Value = 3119.9412
$ws.Range("K138").Value = 5507.6001
$ws.Range("L138").Value = 9359.8236
$ws.Range("M138").Value = -367.6000999999997
$ws.Range("N138").Value = -19639.8236
$ws.Range("H141").Value = 8437.25
$ws.Range("I141").Value = 8205.700000000001
$ws.Range("K141").Value = 24617.1
$ws.Range("M141").Value = -19437.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 499
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 499
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 499
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -849
$ws.Range("H55").Value = 20933.334
$ws.Range("J55").Value = 29900
$ws.Range("L55").Value = 29900
$ws.Range("N55").Value = -30530
$ws.Range("H61").Value = 7405.3706
$ws.Range("J61").Value = 4999.5
$ws.Range("L61").Value = 4999.5
$ws.Range("N61").Value = -5423.5
$ws.Range("H97").Value = 685.53845
$ws.Range("I97").Value = 628.04346
$ws.Range("J97").Value = 1126.3334
$ws.Range("K97").Value = 628.04346
$ws.Range("L97").Value = 1126.3334
$ws.Range("M97").Value = -132.04346
$ws.Range("N97").Value = -2118.3334
$ws.Range("H122").Value = 1918.8889
$ws.Range("I122").Value = 1918.8889
$ws.Range("K122").Value = 5756.6667
$ws.Range("M122").Value = -3306.6667
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080
$ws.Range("H136").Value = 7405.3706
$ws.Range("J136").Value = 4999.5
$ws.Range("L136").Value = 14998.5
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2270.1667
$ws.Range("I107").Value = 2634.2
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 2634.2
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = -714.1999999999998
$ws.Range("N107").Value = -4290
$ws.Range("H134").Value = 3200.6848
$ws.Range("I134").Value = 2855.4575
$ws.Range("J134").Value = 4655.5713
$ws.Range("K134").Value = 8566.372499999999
$ws.Range("L134").Value = 13966.7139
$ws.Range("M134").Value = -6031.372499999999
$ws.Range("N134").Value = -19036.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11870
$ws.Range("I6").Value = 21169
$ws.Range("J6").Value = 246.25
$ws.Range("K6").Value = 21169
$ws.Range("L6").Value = 246.25
$ws.Range("M6").Value = -21056
$ws.Range("N6").Value = -472.25
$ws.Range("H31").Value = 5519.4736
$ws.Range("J31").Value = 10965.6
$ws.Range("L31").Value = 10965.6
$ws.Range("N31").Value = -11555.6
$ws.Range("H34").Value = 5519.4736
$ws.Range("J34").Value = 10965.6
$ws.Range("L34").Value = 10965.6
$ws.Range("N34").Value = -11369.6
$ws.Range("H74").Value = 37140
$ws.Range("J74").Value = 37140
$ws.Range("L74").Value = 37140
$ws.Range("N74").Value = -38888
$ws.Range("H77").Value = 37140
$ws.Range("J77").Value = 37140
$ws.Range("L77").Value = 111420
$ws.Range("N77").Value = -120156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17040956
$ws.Range("I4").Value = 13724687
$ws.Range("J4").Value = 25000000
$ws.Range("K4").Value = 41174061
$ws.Range("L4").Value = 75000000
$ws.Range("M4").Value = -41173949
$ws.Range("N4").Value = -75000224
$ws.Range("H12").Value = 474.64706
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 722.9091
$ws.Range("K12").Value = 58.5
$ws.Range("L12").Value = 2168.7273
$ws.Range("M12").Value = 114.5
$ws.Range("N12").Value = -2514.7273
$ws.Range("H22").Value = 1318.2
$ws.Range("I22").Value = 351.66666
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 1054.99998
$ws.Range("L22").Value = 4350
$ws.Range("M22").Value = -885.9999800000001
$ws.Range("N22").Value = -4688
$ws.Range("H27").Value = 1318.2
$ws.Range("I27").Value = 351.66666
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 1054.99998
$ws.Range("L27").Value = 4350
$ws.Range("M27").Value = -952.9999800000001
$ws.Range("N27").Value = -4554
$ws.Range("H46").Value = 561.1111
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 600
$ws.Range("K46").Value = 1275
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -1184
$ws.Range("N46").Value = -1982
$ws.Range("H52").Value = 2110.3333
$ws.Range("J52").Value = 2110.3333
$ws.Range("L52").Value = 6330.999899999999
$ws.Range("N52").Value = -6862.999899999999
$ws.Range("H107").Value = 280.4
$ws.Range("I107").Value = 228.66667
$ws.Range("J107").Value = 302.57144
$ws.Range("K107").Value = 686.00001
$ws.Range("L107").Value = 907.71432
$ws.Range("M107").Value = 1233.99999
$ws.Range("N107").Value = -4747.71432
$ws.Range("H128").Value = 1129998
$ws.Range("I128").Value = 1129998
$ws.Range("K128").Value = 3389994
$ws.Range("M128").Value = -3385014
$ws.Range("H132").Value = 1402.875
$ws.Range("J132").Value = 1702.3125
$ws.Range("L132").Value = 15320.8125
$ws.Range("N132").Value = -20380.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 18582.666
$ws.Range("J46").Value = 22648.5
$ws.Range("L46").Value = 22648.5
$ws.Range("N46").Value = -22960.5
$ws.Range("H102").Value = 6518.385
$ws.Range("I102").Value = 4523.5557
$ws.Range("K102").Value = 4523.5557
$ws.Range("M102").Value = -2901.5557
$ws.Range("H126").Value = 3511
$ws.Range("J126").Value = 5071.2856
$ws.Range("L126").Value = 15213.8568
$ws.Range("N126").Value = -20153.8568
$ws.Range("H132").Value = 4519.4165
$ws.Range("I132").Value = 4839.2354
$ws.Range("K132").Value = 14517.7062
$ws.Range("M132").Value = -11987.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3943.0588
$ws.Range("I40").Value = 4008.4
$ws.Range("K40").Value = 4008.4
$ws.Range("M40").Value = -3872.4
$ws.Range("H46").Value = 2377.3157
$ws.Range("I46").Value = 1894.0834
$ws.Range("K46").Value = 1894.0834
$ws.Range("M46").Value = -1706.0834
$ws.Range("H93").Value = 10730.5
$ws.Range("I93").Value = 1219.5
$ws.Range("J93").Value = 27374.75
$ws.Range("K93").Value = 1219.5
$ws.Range("L93").Value = 27374.75
$ws.Range("M93").Value = 28.5
$ws.Range("N93").Value = -29870.75
$ws.Range("H136").Value = 47623876
$ws.Range("I136").Value = 25004822
$ws.Range("K136").Value = 75014466
$ws.Range("M136").Value = -75011916
$ws.Range("H139").Value = 40000
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -34860
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 6982.5625
$ws.Range("I15").Value = 6983.846
$ws.Range("J15").Value = 6977
$ws.Range("K15").Value = 6983.846
$ws.Range("L15").Value = 6977
$ws.Range("M15").Value = -6695.846
$ws.Range("N15").Value = -7553
$ws.Range("H96").Value = 3999.6667
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3999.6667
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3999.6667
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -6745.6667
$ws.Range("H107").Value = 1306.625
$ws.Range("I107").Value = 1368.8334
$ws.Range("J107").Value = 1120
$ws.Range("K107").Value = 4106.5002
$ws.Range("L107").Value = 3360
$ws.Range("M107").Value = -2186.5002
$ws.Range("N107").Value = -7200
$ws.Range("H113").Value = 593.9545000000001
$ws.Range("I113").Value = 374.08334
$ws.Range("K113").Value = 1122.25002
$ws.Range("M113").Value = 1047.74998
$ws.Range("H122").Value = 3509.3914
$ws.Range("I122").Value = 2294.75
$ws.Range("K122").Value = 6884.25
$ws.Range("M122").Value = -4434.25
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
